$wb = $excel.ActiveWorkbook

# --- "Diff" sheet: update header cells D1 and F1 with new shared strings ---
$wsDiff = $wb.Worksheets.Item("Diff")
$wsDiff.Range("D1").Value = "col3"
$wsDiff.Range("F1").Value = "col5"

# --- Switch the active/selected tab from "Current" to "Diff" ---
$wsDiff.Activate()

# --- Update the selection on the "Diff" sheet ---
[void]$wsDiff.Range("F4").Select()
